# Appends 11 new training rows (J-3, 2025-12-03) to Feuil1, rows 951-961,
# mirroring the source workbook diff (new GPS session data + matching
# "01:xx:xx" duration strings added to the shared-string table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstNewRow = 951
$lastNewRow = 961

# Seed every new row from row 949s formatting (it already carries the
# correct date style on column B and the centered "J-x" style on column D),
# so the new cells inherit the same cellXfs entries instead of minting new ones.
for ($r = $firstNewRow; $r -le $lastNewRow; $r++) {
    $ws.Range("A949:V949").Copy($ws.Range("A" + $r + ":V" + $r))
}

# Row 951s player name (Ryad Kralladi) is styled with the alternate font used
# elsewhere for this player (e.g. E245) - replicate that one-off style too.
$ws.Range("E245").Copy($ws.Range("E951"))

# Row 951
$ws.Range("A951").Value = "Entrainement"
$ws.Range("B951").Value = 45994
$ws.Range("C951").Value = "Global"
$ws.Range("D951").Value = "J-3"
$ws.Range("E951").Value = "Ryad Kralladi"
$ws.Range("F951").Value = "center midfield"
$ws.Range("G951").Value = "01:43:27"
$ws.Range("H951").Value = 6.35
$ws.Range("I951").Value = 0.67
$ws.Range("J951").Value = 5.67
$ws.Range("K951").Value = 0.46
$ws.Range("L951").Value = 0.18
$ws.Range("M951").Value = 0.04
$ws.Range("N951").Value = 0
$ws.Range("O951").Value = 2
$ws.Range("P951").Value = 3.52
$ws.Range("Q951").Value = 29.59
$ws.Range("R951").Value = 4.1
$ws.Range("S951").Value = 22
$ws.Range("T951").Value = 1
$ws.Range("U951").Value = 13
$ws.Range("V951").Value = 2

# Row 952
$ws.Range("A952").Value = "Entrainement"
$ws.Range("B952").Value = 45994
$ws.Range("C952").Value = "Global"
$ws.Range("D952").Value = "J-3"
$ws.Range("E952").Value = "Amir Etien"
$ws.Range("F952").Value = "right forward"
$ws.Range("G952").Value = "01:43:50"
$ws.Range("H952").Value = 7.04
$ws.Range("I952").Value = 0.92
$ws.Range("J952").Value = 6.1
$ws.Range("K952").Value = 0.57
$ws.Range("L952").Value = 0.2
$ws.Range("M952").Value = 0.1
$ws.Range("N952").Value = 0.06
$ws.Range("O952").Value = 11
$ws.Range("P952").Value = 3.73
$ws.Range("Q952").Value = 35.02
$ws.Range("R952").Value = 6.37
$ws.Range("S952").Value = 64
$ws.Range("T952").Value = 27
$ws.Range("U952").Value = 52
$ws.Range("V952").Value = 14

# Row 953
$ws.Range("A953").Value = "Entrainement"
$ws.Range("B953").Value = 45994
$ws.Range("C953").Value = "Global"
$ws.Range("D953").Value = "J-3"
$ws.Range("E953").Value = "Kamal Bafounta"
$ws.Range("F953").Value = "center midfield"
$ws.Range("G953").Value = "01:42:56"
$ws.Range("H953").Value = 7.17
$ws.Range("I953").Value = 0.92
$ws.Range("J953").Value = 6.24
$ws.Range("K953").Value = 0.69
$ws.Range("L953").Value = 0.15
$ws.Range("M953").Value = 0.07
$ws.Range("N953").Value = 0.02
$ws.Range("O953").Value = 7
$ws.Range("P953").Value = 4.13
$ws.Range("Q953").Value = 31.13
$ws.Range("R953").Value = 4.31
$ws.Range("S953").Value = 24
$ws.Range("T953").Value = 4
$ws.Range("U953").Value = 17
$ws.Range("V953").Value = 3

# Row 954
$ws.Range("A954").Value = "Entrainement"
$ws.Range("B954").Value = 45994
$ws.Range("C954").Value = "Global"
$ws.Range("D954").Value = "J-3"
$ws.Range("E954").Value = "Naim Ighbane"
$ws.Range("F954").Value = "center back"
$ws.Range("G954").Value = "01:43:35"
$ws.Range("H954").Value = 7.28
$ws.Range("I954").Value = 0.51
$ws.Range("J954").Value = 6.76
$ws.Range("K954").Value = 0.48
$ws.Range("L954").Value = 0.04
$ws.Range("M954").Value = 0
$ws.Range("N954").Value = 0
$ws.Range("O954").Value = 0
$ws.Range("P954").Value = 3.8
$ws.Range("Q954").Value = 24.96
$ws.Range("R954").Value = 5.19
$ws.Range("S954").Value = 50
$ws.Range("T954").Value = 12
$ws.Range("U954").Value = 46
$ws.Range("V954").Value = 9

# Row 955
$ws.Range("A955").Value = "Entrainement"
$ws.Range("B955").Value = 45994
$ws.Range("C955").Value = "Global"
$ws.Range("D955").Value = "J-3"
$ws.Range("E955").Value = "Malik Boussaid"
$ws.Range("F955").Value = "right back"
$ws.Range("G955").Value = "01:43:42"
$ws.Range("H955").Value = 7.98
$ws.Range("I955").Value = 1.13
$ws.Range("J955").Value = 6.84
$ws.Range("K955").Value = 0.52
$ws.Range("L955").Value = 0.38
$ws.Range("M955").Value = 0.18
$ws.Range("N955").Value = 0.06
$ws.Range("O955").Value = 13
$ws.Range("P955").Value = 4.2
$ws.Range("Q955").Value = 34.48
$ws.Range("R955").Value = 5.17
$ws.Range("S955").Value = 66
$ws.Range("T955").Value = 8
$ws.Range("U955").Value = 55
$ws.Range("V955").Value = 9

# Row 956
$ws.Range("A956").Value = "Entrainement"
$ws.Range("B956").Value = 45994
$ws.Range("C956").Value = "Global"
$ws.Range("D956").Value = "J-3"
$ws.Range("E956").Value = "Mattheo Haon"
$ws.Range("F956").Value = "right back"
$ws.Range("G956").Value = "01:43:50"
$ws.Range("H956").Value = 8.04
$ws.Range("I956").Value = 1.06
$ws.Range("J956").Value = 6.97
$ws.Range("K956").Value = 0.55
$ws.Range("L956").Value = 0.31
$ws.Range("M956").Value = 0.19
$ws.Range("N956").Value = 0.03
$ws.Range("O956").Value = 10
$ws.Range("P956").Value = 4.6
$ws.Range("Q956").Value = 30.85
$ws.Range("R956").Value = 4.45
$ws.Range("S956").Value = 32
$ws.Range("T956").Value = 3
$ws.Range("U956").Value = 33
$ws.Range("V956").Value = 4

# Row 957
$ws.Range("A957").Value = "Entrainement"
$ws.Range("B957").Value = 45994
$ws.Range("C957").Value = "Global"
$ws.Range("D957").Value = "J-3"
$ws.Range("E957").Value = "Emmanuel Valey"
$ws.Range("F957").Value = "left forward"
$ws.Range("G957").Value = "01:42:35"
$ws.Range("H957").Value = 8.05
$ws.Range("I957").Value = 1.35
$ws.Range("J957").Value = 6.67
$ws.Range("K957").Value = 0.79
$ws.Range("L957").Value = 0.4
$ws.Range("M957").Value = 0.17
$ws.Range("N957").Value = 0.02
$ws.Range("O957").Value = 15
$ws.Range("P957").Value = 4.18
$ws.Range("Q957").Value = 31.28
$ws.Range("R957").Value = 4.77
$ws.Range("S957").Value = 51
$ws.Range("T957").Value = 9
$ws.Range("U957").Value = 46
$ws.Range("V957").Value = 11

# Row 958
$ws.Range("A958").Value = "Entrainement"
$ws.Range("B958").Value = 45994
$ws.Range("C958").Value = "Global"
$ws.Range("D958").Value = "J-3"
$ws.Range("E958").Value = "Yoan Zouma"
$ws.Range("F958").Value = "center back"
$ws.Range("G958").Value = "01:41:57"
$ws.Range("H958").Value = 6.14
$ws.Range("I958").Value = 0.41
$ws.Range("J958").Value = 5.73
$ws.Range("K958").Value = 0.34
$ws.Range("L958").Value = 0.07
$ws.Range("M958").Value = 0
$ws.Range("N958").Value = 0
$ws.Range("O958").Value = 0
$ws.Range("P958").Value = 3.5
$ws.Range("Q958").Value = 22.49
$ws.Range("R958").Value = 4.4
$ws.Range("S958").Value = 27
$ws.Range("T958").Value = 5
$ws.Range("U958").Value = 13
$ws.Range("V958").Value = 1

# Row 959
$ws.Range("A959").Value = "Entrainement"
$ws.Range("B959").Value = 45994
$ws.Range("C959").Value = "Global"
$ws.Range("D959").Value = "J-3"
$ws.Range("E959").Value = "Jeremie Laurent"
$ws.Range("F959").Value = "left forward"
$ws.Range("G959").Value = "01:19:52"
$ws.Range("H959").Value = 6.14
$ws.Range("I959").Value = 0.6
$ws.Range("J959").Value = 5.53
$ws.Range("K959").Value = 0.49
$ws.Range("L959").Value = 0.12
$ws.Range("M959").Value = 0
$ws.Range("N959").Value = 0
$ws.Range("O959").Value = 0
$ws.Range("P959").Value = 4.55
$ws.Range("Q959").Value = 24.65
$ws.Range("R959").Value = 4.82
$ws.Range("S959").Value = 29
$ws.Range("T959").Value = 5
$ws.Range("U959").Value = 15
$ws.Range("V959").Value = 1

# Row 960
$ws.Range("A960").Value = "Entrainement"
$ws.Range("B960").Value = 45994
$ws.Range("C960").Value = "Global"
$ws.Range("D960").Value = "J-3"
$ws.Range("E960").Value = "Sofiane Belle"
$ws.Range("F960").Value = "left forward"
$ws.Range("G960").Value = "01:43:19"
$ws.Range("H960").Value = 7.12
$ws.Range("I960").Value = 1.1
$ws.Range("J960").Value = 6.01
$ws.Range("K960").Value = 0.67
$ws.Range("L960").Value = 0.35
$ws.Range("M960").Value = 0.08
$ws.Range("N960").Value = 0.01
$ws.Range("O960").Value = 8
$ws.Range("P960").Value = 4.02
$ws.Range("Q960").Value = 30.63
$ws.Range("R960").Value = 4.29
$ws.Range("S960").Value = 22
$ws.Range("T960").Value = 1
$ws.Range("U960").Value = 17
$ws.Range("V960").Value = 6

# Row 961
$ws.Range("A961").Value = "Entrainement"
$ws.Range("B961").Value = 45994
$ws.Range("C961").Value = "Global"
$ws.Range("D961").Value = "J-3"
$ws.Range("E961").Value = "Karim Belmahi"
$ws.Range("F961").Value = "left forward"
$ws.Range("G961").Value = "01:43:27"
$ws.Range("H961").Value = 7.06
$ws.Range("I961").Value = 0.74
$ws.Range("J961").Value = 6.31
$ws.Range("K961").Value = 0.49
$ws.Range("L961").Value = 0.19
$ws.Range("M961").Value = 0.08
$ws.Range("N961").Value = 0
$ws.Range("O961").Value = 6
$ws.Range("P961").Value = 3.53
$ws.Range("Q961").Value = 29.4
$ws.Range("R961").Value = 4.75
$ws.Range("S961").Value = 35
$ws.Range("T961").Value = 3
$ws.Range("U961").Value = 20
$ws.Range("V961").Value = 8
